$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: insert 3 new empty paragraphs (style KeyHeadDetails) right after
# "Our Ref: <Primary Reference Number>" and before the paragraph that holds
# the 4 tab characters.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Our Ref: <Primary Reference Number>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $ourRefPara = $rng.Paragraphs(1)
    $insertPoint = $d.Range($ourRefPara.Range.End, $ourRefPara.Range.End)
    $xmlFragParas = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="KeyHeadDetails"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="KeyHeadDetails"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="KeyHeadDetails"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insertPoint.InsertXML($xmlFragParas)
}

# ---------------------------------------------------------------------------
# Change 2: "NATIONAL PLANNING POLICY FRAMEWORK 2021" -> "... 2023"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("NATIONAL PLANNING POLICY FRAMEWORK 2021", $true, $false, $false, $false, $false, $true, 1, $false, "NATIONAL PLANNING POLICY FRAMEWORK 2023", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 3: <Proposal Description> paragraph gains italic (i + iCs) on top
# of the existing bold, both on the run and on the paragraph mark.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("<Proposal Description>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $pdPara = $rng.Paragraphs(1)
    $xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="5245"/></w:tabs><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/><w:i/><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs="Arial"/><w:b/><w:i/><w:iCs/></w:rPr><w:t>&lt;Proposal Description&gt;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $pdPara.Range.InsertXML($xmlFrag)
}

# ---------------------------------------------------------------------------
# Change 4: renumber the NPPF paragraph references / years.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("NPPF paragraph 194 says", $true, $false, $false, $false, $false, $true, 1, $false, "NPPF paragraph 200 says", 2) | Out-Null
$d.Content.Find.Execute("NPPF paragraphs 199 - 202 place great weight", $true, $false, $false, $false, $false, $true, 1, $false, "NPPF paragraphs 205 - 208 place great weight", 2) | Out-Null
$d.Content.Find.Execute("(NPPF paragraph 203)", $true, $false, $false, $false, $false, $true, 1, $false, "(NPPF paragraph 209)", 2) | Out-Null
$d.Content.Find.Execute("NPPF paragraphs 190 and 197 and London Plan", $true, $false, $false, $false, $false, $true, 1, $false, "NPPF paragraphs 195 and 203 and London Plan", 2) | Out-Null
$d.Content.Find.Execute("paragraph 205 of the NPPF says", $true, $false, $false, $false, $false, $true, 1, $false, "paragraph 211 of the NPPF says", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 5: the <Casework Officer> paragraph right after "Yours sincerely"
# becomes bold (b + bCs) on both the run and the paragraph mark.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Yours sincerely", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $yoursSincerelyPara = $rng.Paragraphs(1)
    $caseworkPara = $yoursSincerelyPara.Next().Next()
    $xmlFrag2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:tabs><w:tab w:val="left" w:pos="5245"/></w:tabs><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>&lt;Casework Officer&gt;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $caseworkPara.Range.InsertXML($xmlFrag2)
}
